$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 190, shifting existing rows 190-196 down to 191-197
$ws.Rows.Item(190).EntireRow.Insert()

# Fill in the new row 190 with the new weekly price record
$ws.Cells.Item(190, 1).Value = 11
$ws.Cells.Item(190, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(190, 3).Value = "Bíobío"
$ws.Cells.Item(190, 4).Value = 44714
$ws.Cells.Item(190, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(190, 5).Value = 8
$ws.Cells.Item(190, 6).Value = 100112040
$ws.Cells.Item(190, 7).Value = "Cilantro"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 160
$ws.Cells.Item(190, 11).Value = 5000
$ws.Cells.Item(190, 12).Value = 5500
$ws.Cells.Item(190, 13).Value = 5250
$ws.Cells.Item(190, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(190, 15).Value = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value = 146
$ws.Cells.Item(190, 17).Value = 36
$ws.Cells.Item(190, 18).Value = "Hortaliza"
